# Generate Report for Archive
#
# The localization status report is being refreshed: every row that was
# previously marked "Ready for handoff" has moved on to "In Translation".
# That text lives in a shared string that is referenced from the per-locale
# "Status" column (zh-cn / de-de sheets) and is mirrored onto the "Overview"
# sheet's zh-cn / de-de columns. Once the text shrinks, those same "Status"
# columns are re-sized to fit the new (shorter) content.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status mirrors ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"
$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"
$ws.Range("E1:F1").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (col C) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"
$ws.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: Status column (col C) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"
$ws.Range("C1").ColumnWidth = 12.5
